# Update the "dSF" column (F) values for several rows, per repull/push of data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F3").Value = 2
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = 3
$ws.Range("F10").Value = 5
